$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "COD" column (K), matching the formatting
# already used by the adjacent empty cells (thin box border, no fill).
$ws.Range("J4").Copy()
$ws.Range("K1:K13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header cell for the new column.
$ws.Range("K1").Value = "COD"

# Sample data value under the new header (second data row).
$ws.Range("K2").Value = 4

# Reflect the new selection state left by the edit.
[void]$ws.Range("K1:K13").Select()
